# Scheduled market-data refresh: update currentAveragePrice* / Leve*Price / Leve*Profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4717.727
$ws.Range("I64").Value = 4262.5
$ws.Range("K64").Value = 4262.5
$ws.Range("M64").Value = -4014.5
# Row 67
$ws.Range("H67").Value = 4717.727
$ws.Range("I67").Value = 4262.5
$ws.Range("K67").Value = 4262.5
$ws.Range("M67").Value = -3404.5
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 80
$ws.Range("H80").Value = 1796.9231
$ws.Range("I80").Value = 1008
$ws.Range("J80").Value = 3059.2
$ws.Range("K80").Value = 3024
$ws.Range("L80").Value = 9177.599999999999
$ws.Range("M80").Value = -2026
$ws.Range("N80").Value = -11173.6
# Row 83
$ws.Range("H83").Value = 1796.9231
$ws.Range("I83").Value = 1008
$ws.Range("J83").Value = 3059.2
$ws.Range("K83").Value = 9072
$ws.Range("L83").Value = 27532.8
$ws.Range("M83").Value = -4080
$ws.Range("N83").Value = -37516.8
# Row 106
$ws.Range("H106").Value = 4885.2
$ws.Range("I106").Value = 5069.857
$ws.Range("K106").Value = 5069.857
$ws.Range("M106").Value = -4438.857
# Row 137
$ws.Range("H137").Value = 4648.6387
$ws.Range("I137").Value = 1926.7333
$ws.Range("J137").Value = 18258.166
$ws.Range("K137").Value = 5780.199900000001
$ws.Range("L137").Value = 54774.49800000001
$ws.Range("M137").Value = -3230.199900000001
$ws.Range("N137").Value = -59874.49800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3342.5715
$ws.Range("J63").Value = 5266.3335
$ws.Range("L63").Value = 5266.3335
$ws.Range("N63").Value = -6638.3335
# Row 66
$ws.Range("H66").Value = 3342.5715
$ws.Range("J66").Value = 5266.3335
$ws.Range("L66").Value = 26331.6675
$ws.Range("N66").Value = -33195.6675
# Row 97
$ws.Range("H97").Value = 948.9
$ws.Range("I97").Value = 830.4211
$ws.Range("J97").Value = 3200
$ws.Range("K97").Value = 830.4211
$ws.Range("L97").Value = 3200
$ws.Range("M97").Value = -334.4211
$ws.Range("N97").Value = -4192

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 26437
$ws.Range("I64").Value = 100000
$ws.Range("J64").Value = 1916
$ws.Range("K64").Value = 100000
$ws.Range("L64").Value = 1916
$ws.Range("M64").Value = -99775
$ws.Range("N64").Value = -2366
# Row 67
$ws.Range("H67").Value = 26437
$ws.Range("I67").Value = 100000
$ws.Range("J67").Value = 1916
$ws.Range("K67").Value = 100000
$ws.Range("L67").Value = 1916
$ws.Range("M67").Value = -99220
$ws.Range("N67").Value = -3476
# Row 96
$ws.Range("H96").Value = 35080.2
$ws.Range("I96").Value = 21800.334
$ws.Range("K96").Value = 21800.334
$ws.Range("M96").Value = -19054.334
# Row 99
$ws.Range("H99").Value = 16170.571
$ws.Range("I99").Value = 16170.571
$ws.Range("K99").Value = 16170.571
$ws.Range("M99").Value = -14672.571

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 4791.1665
$ws.Range("I62").Value = 4687
$ws.Range("K62").Value = 4687
$ws.Range("M62").Value = -4063
# Row 65
$ws.Range("H65").Value = 4791.1665
$ws.Range("I65").Value = 4687
$ws.Range("K65").Value = 23435
$ws.Range("M65").Value = -20315
# Row 99
$ws.Range("H99").Value = 31224.572
$ws.Range("I99").Value = 52089.5
$ws.Range("K99").Value = 52089.5
$ws.Range("M99").Value = -50591.5
# Row 126
$ws.Range("H126").Value = 31224.572
$ws.Range("I126").Value = 52089.5
$ws.Range("K126").Value = 156268.5
$ws.Range("M126").Value = -153798.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 4817020
$ws.Range("I4").Value = 9286442
$ws.Range("K4").Value = 27859326
$ws.Range("M4").Value = -27859214
# Row 32
$ws.Range("H32").Value = 10000000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10000000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 30000000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -30000566
# Row 34
$ws.Range("H34").Value = 1762.5883
$ws.Range("J34").Value = 1894.9231
$ws.Range("L34").Value = 5684.7693
$ws.Range("N34").Value = -5852.7693
# Row 46
$ws.Range("H46").Value = 455506.2
$ws.Range("J46").Value = 2502124.2
$ws.Range("L46").Value = 7506372.600000001
$ws.Range("N46").Value = -7506554.600000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5043.3335
$ws.Range("J80").Value = 4583.3335
$ws.Range("L80").Value = 4583.3335
$ws.Range("N80").Value = -6579.3335
# Row 83
$ws.Range("H83").Value = 5043.3335
$ws.Range("J83").Value = 4583.3335
$ws.Range("L83").Value = 22916.6675
$ws.Range("N83").Value = -32900.6675
# Row 111
$ws.Range("H111").Value = 47500
$ws.Range("J111").Value = 47500
$ws.Range("L111").Value = 47500
$ws.Range("N111").Value = -53634

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 9455.368
$ws.Range("I68").Value = 8691.416999999999
$ws.Range("J68").Value = 10765
$ws.Range("K68").Value = 8691.416999999999
$ws.Range("L68").Value = 10765
$ws.Range("M68").Value = -7942.416999999999
$ws.Range("N68").Value = -12263
# Row 71
$ws.Range("H71").Value = 9455.368
$ws.Range("I71").Value = 8691.416999999999
$ws.Range("J71").Value = 10765
$ws.Range("K71").Value = 43457.085
$ws.Range("L71").Value = 53825
$ws.Range("M71").Value = -39713.085
$ws.Range("N71").Value = -61313
# Row 82
$ws.Range("H82").Value = 2190.4167
$ws.Range("I82").Value = 2223.125
$ws.Range("J82").Value = 2125
$ws.Range("K82").Value = 2223.125
$ws.Range("L82").Value = 2125
$ws.Range("M82").Value = -1862.125
$ws.Range("N82").Value = -2847
# Row 85
$ws.Range("H85").Value = 2190.4167
$ws.Range("I85").Value = 2223.125
$ws.Range("J85").Value = 2125
$ws.Range("K85").Value = 2223.125
$ws.Range("L85").Value = 2125
$ws.Range("M85").Value = -975.125
$ws.Range("N85").Value = -4621
# Row 93
$ws.Range("H93").Value = 3286.4546
$ws.Range("I93").Value = 1574.1666
$ws.Range("J93").Value = 5341.2
$ws.Range("K93").Value = 1574.1666
$ws.Range("L93").Value = 5341.2
$ws.Range("M93").Value = -326.1666
$ws.Range("N93").Value = -7837.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 14744.167
$ws.Range("I62").Value = 5165.778
$ws.Range("K62").Value = 5165.778
$ws.Range("M62").Value = -4541.778
# Row 65
$ws.Range("H65").Value = 14744.167
$ws.Range("I65").Value = 5165.778
$ws.Range("K65").Value = 25828.89
$ws.Range("M65").Value = -22708.89
# Row 119
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
